$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Julio de 2020 a las 18:32"

# --- Swap country names for rows whose rank changed due to updated case counts ---
# Barein (row 50) <-> Israel (row 51)
$ws.Cells.Item(50,1).Value = "Israel"
$ws.Cells.Item(51,1).Value = "Barein"
# Zimbabue (row 152) <-> Montenegro (row 153)
$ws.Cells.Item(152,1).Value = "Montenegro"
$ws.Cells.Item(153,1).Value = "Zimbabue"
# Dominica (row 205) <-> Fiyi (row 206)
$ws.Cells.Item(205,1).Value = "Fiyi"
$ws.Cells.Item(206,1).Value = "Dominica"
# Islas Malvinas (row 209) <-> Groenlandia (row 210)
$ws.Cells.Item(209,1).Value = "Groenlandia"
$ws.Cells.Item(210,1).Value = "Islas Malvinas"

# --- Update statistic values (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Cells.Item(4,2).Value = 2857996
$ws.Cells.Item(4,3).Value = 22312
$ws.Cells.Item(4,4).Value = 1196353
$ws.Cells.Item(4,5).Value = 1529949
$ws.Cells.Item(4,7).Value = 209
$ws.Cells.Item(4,8).Value = 131694

$ws.Cells.Item(7,2).Value = 644404
$ws.Cells.Item(7,3).Value = 17236
$ws.Cells.Item(7,4).Value = 390252
$ws.Cells.Item(7,5).Value = 235555
$ws.Cells.Item(7,7).Value = 372
$ws.Cells.Item(7,8).Value = 18597

$ws.Cells.Item(8,2).Value = 297625
$ws.Cells.Item(8,3).Value = 442
$ws.Cells.Item(8,7).Value = 17
$ws.Cells.Item(8,8).Value = 28385

$ws.Cells.Item(10,2).Value = 288089
$ws.Cells.Item(10,3).Value = 3548
$ws.Cells.Item(10,4).Value = 253343
$ws.Cells.Item(10,5).Value = 28695
$ws.Cells.Item(10,7).Value = 131
$ws.Cells.Item(10,8).Value = 6051

$ws.Cells.Item(11,2).Value = 284276
$ws.Cells.Item(11,3).Value = 519
$ws.Cells.Item(11,7).Value = 136
$ws.Cells.Item(11,8).Value = 44131

$ws.Cells.Item(12,2).Value = 241184
$ws.Cells.Item(12,3).Value = 223
$ws.Cells.Item(12,4).Value = 191467
$ws.Cells.Item(12,5).Value = 14884
$ws.Cells.Item(12,7).Value = 15
$ws.Cells.Item(12,8).Value = 34833

$ws.Cells.Item(23,2).Value = 104936
$ws.Cells.Item(23,3).Value = 164
$ws.Cells.Item(23,4).Value = 68526
$ws.Cells.Item(23,5).Value = 27766
$ws.Cells.Item(23,7).Value = 2
$ws.Cells.Item(23,8).Value = 8644

$ws.Cells.Item(38,4).Value = 39769
$ws.Cells.Item(38,5).Value = 4684

$ws.Cells.Item(50,2).Value = 28055
$ws.Cells.Item(50,3).Value = 1008
$ws.Cells.Item(50,4).Value = 17669
$ws.Cells.Item(50,5).Value = 10060
$ws.Cells.Item(50,7).Value = 2
$ws.Cells.Item(50,8).Value = 326

$ws.Cells.Item(51,2).Value = 27837
$ws.Cells.Item(51,3).Value = 0
$ws.Cells.Item(51,4).Value = 22583
$ws.Cells.Item(51,5).Value = 5159
$ws.Cells.Item(51,7).Value = 1
$ws.Cells.Item(51,8).Value = 95

$ws.Cells.Item(69,2).Value = 12248
$ws.Cells.Item(69,3).Value = 70
$ws.Cells.Item(69,5).Value = 4073

$ws.Cells.Item(72,4).Value = 6251
$ws.Cells.Item(72,5).Value = 3047

$ws.Cells.Item(98,2).Value = 3486
$ws.Cells.Item(98,3).Value = 28
$ws.Cells.Item(98,5).Value = 1920

$ws.Cells.Item(112,2).Value = 2069
$ws.Cells.Item(112,3).Value = 3
$ws.Cells.Item(112,5).Value = 195

$ws.Cells.Item(115,4).Value = 1832
$ws.Cells.Item(115,5).Value = 13

$ws.Cells.Item(126,2).Value = 1382
$ws.Cells.Item(126,3).Value = 81
$ws.Cells.Item(126,5).Value = 738

$ws.Cells.Item(130,2).Value = 1181
$ws.Cells.Item(130,3).Value = 3
$ws.Cells.Item(130,4).Value = 1045
$ws.Cells.Item(130,5).Value = 86

$ws.Cells.Item(131,2).Value = 1147
$ws.Cells.Item(131,3).Value = 11
$ws.Cells.Item(131,4).Value = 897
$ws.Cells.Item(131,5).Value = 241

$ws.Cells.Item(139,2).Value = 939
$ws.Cells.Item(139,3).Value = 21
$ws.Cells.Item(139,5).Value = 684

$ws.Cells.Item(152,2).Value = 663
$ws.Cells.Item(152,3).Value = 47
$ws.Cells.Item(152,4).Value = 315
$ws.Cells.Item(152,5).Value = 335
$ws.Cells.Item(152,7).Value = 1
$ws.Cells.Item(152,8).Value = 13

$ws.Cells.Item(153,2).Value = 617
$ws.Cells.Item(153,4).Value = 173
$ws.Cells.Item(153,5).Value = 437
$ws.Cells.Item(153,8).Value = 7
